# Apply the cryptos list refresh (GitHub Actions data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.343.73"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "2.035.40"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.70"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.655"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.80"
$ws.Range("E8").Value = "  -8.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.67"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.358"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0736"
$ws.Range("E11").Value = "  -4.78%  "
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.923"
$ws.Range("E13").Value = "  +6.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.39"
$ws.Range("E14").Value = "  -5.06%  "
$ws.Range("D15").Value = "2.333.88"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "2.037.18"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "36.280.45"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.76"
$ws.Range("E19").Value = "  -6.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.88"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.20"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.10"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.28"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("E28").Value = "  -11.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.61"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.90"
$ws.Range("E32").Value = "  -10.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0584"
$ws.Range("E33").Value = "  -4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  -7.16%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0862"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.18"
$ws.Range("E38").Value = "  -5.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.90"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.21"
$ws.Range("E40").Value = "  -7.73%  "
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0211"
$ws.Range("E42").Value = "  -5.14%  "
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "91.91"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("D46").Value = "1.367.64"
$ws.Range("E46").Value = "  +5.42%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.46"
$ws.Range("E47").Value = "  -7.86%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.31"
$ws.Range("E48").Value = "  +8.97%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "2.216.98"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  -4.89%  "
